# CSIEM Data Catalogue — add new "WAMSI Westport Marine Science Program"
# (WWMSP3) row to the catalogue table on Sheet1.
#
# The new record is inserted as row 35 (pushing the existing UWA rows
# 35-38 down to 36-39), directly below the existing WWMSP5.1 row (34),
# since both rows describe the same Program/Agency but a different
# Program Code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 35, shifting rows 35:38 down to 36:39.
$ws.Rows("35").Insert()

# Match the row height used by its neighbours (25.5pt, same as every
# other wrapped-text data row in the table).
$ws.Rows("35").RowHeight = 25.5

# Populate the new record.
$ws.Cells.Item(35, 1).Value = "Data"
$ws.Cells.Item(35, 2).Value = "State Programs"
$ws.Cells.Item(35, 3).Value = "Western Australian Marine Science Institution"
$ws.Cells.Item(35, 4).Value = "WAMSI"
$ws.Cells.Item(35, 5).Value = "WAMSI Westport Marine Science Program"
$ws.Cells.Item(35, 6).Value = "WWMSP3"
$ws.Cells.Item(35, 7).Value = "WQ Grab"
$ws.Cells.Item(35, 8).Value = 18
$ws.Cells.Item(35, 9).Value = "Ongoing"
$ws.Cells.Item(35, 10).Value = "Y"

# Reflect the author's final cursor position (selection moved onto the
# newly-shifted J36 cell of the "SCEVO" row while reviewing the edit).
$ws.Activate() | Out-Null
$ws.Range("J36").Select() | Out-Null
